# DxH_ConfigImport_Veniano_HA.xlsx - import feature changes
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename "commonparameterstest" -> "commonparameters"
# ---------------------------------------------------------------------------
$wsCommon = $wb.Worksheets.Item("commonparameterstest")
$wsCommon.Name = "commonparameters"

# ---------------------------------------------------------------------------
# 2. dtreason sheet: asset_id (col M) shifts down by 8 for every data row
# ---------------------------------------------------------------------------
$wsReason = $wb.Worksheets.Item("dtreason")
for ($r = 2; $r -le 222; $r++) {
    $cell = $wsReason.Cells.Item($r, 13)
    $cell.Value = $cell.Value2 - 8
}
$wsReason.Range("M210:M222").Select()

# ---------------------------------------------------------------------------
# 3. shift sheet: asset_id (col R) 127 -> 119
# ---------------------------------------------------------------------------
$wsShift = $wb.Worksheets.Item("shift")
$wsShift.Range("R2").Value = 119
$wsShift.Range("R3").Value = 119
$wsShift.Range("R3").Select()

# ---------------------------------------------------------------------------
# 4. tag sheet: new "site_id" column inserted before asset_id, plus 12 new
#    tag rows (168H-0001 .. 168H-0013, excluding the pre-existing 168H-0012)
# ---------------------------------------------------------------------------
$wsTag = $wb.Worksheets.Item("tag")

# Header row: P1 becomes site_id, new Q1 becomes asset_id
$wsTag.Range("P1").Value = "site_id"
$wsTag.Range("Q1").Value = "asset_id"

$tagRows = @(
    @{A=56; B="168H-0001.Cutter.Count"; D="NULL"; M="2019-12-17 00:00:00.000"; N="SQL"; O="2019-12-17 00:00:00.000"; P=119; Q=101}
    @{A=59; B="168H-0004.Cutter.Count"; D="NULL"; M="2019-12-17 00:00:00.000"; N="SQL"; O="2019-12-17 00:00:00.000"; P=119; Q=104}
    @{A=60; B="168H-0005.Cutter.Count"; D="NULL"; M="2019-12-17 00:00:00.000"; N="SQL"; O="2019-12-17 00:00:00.000"; P=119; Q=105}
    @{A=61; B="168H-0006.Cutter.Count"; D="NULL"; M="2019-12-17 00:00:00.000"; N="SQL"; O="2019-12-17 00:00:00.000"; P=119; Q=106}
    @{A=62; B="168H-0007.Cutter.Count"; D="NULL"; M="2019-12-17 00:00:00.000"; N="SQL"; O="2019-12-17 00:00:00.000"; P=119; Q=107}
    @{A=63; B="168H-0008.Cutter.Count"; D="NULL"; M="2019-12-17 00:00:00.000"; N="SQL"; O="2019-12-17 00:00:00.000"; P=119; Q=108}
    @{A=64; B="168H-0009.Cutter.Count"; D="NULL"; M="2019-12-17 00:00:00.000"; N="SQL"; O="2019-12-17 00:00:00.000"; P=119; Q=109}
    @{A=65; B="168H-0010.Cutter.Count"; D="NULL"; M="2019-12-17 00:00:00.000"; N="SQL"; O="2019-12-17 00:00:00.000"; P=119; Q=110}
    @{A=66; B="168H-0011.Cutter.Count"; D="NULL"; M="2019-12-17 00:00:00.000"; N="SQL"; O="2019-12-17 00:00:00.000"; P=119; Q=111}
    @{A=67; B="168H-0013.Cutter.Count"; D="NULL"; M="2019-12-17 00:00:00.000"; N="SQL"; O="2019-12-17 00:00:00.000"; P=119; Q=113}
    @{A=57; B="168H-0002.Cutter.Count"; D="NULL"; M="2019-12-17 00:00:00.000"; N="SQL"; O="2019-12-17 00:00:00.000"; P=119; Q=102}
    @{A=58; B="168H-0003.Cutter.Count"; D="NULL"; M="2019-12-17 00:00:00.000"; N="SQL"; O="2019-12-17 00:00:00.000"; P=119; Q=103}
)

$r = 2
foreach ($row in $tagRows) {
    $wsTag.Cells.Item($r, 1).Value = $row.A
    $wsTag.Cells.Item($r, 2).Value = $row.B
    $wsTag.Cells.Item($r, 3).Value = $row.B
    $wsTag.Cells.Item($r, 4).Value = $row.D
    $wsTag.Cells.Item($r, 5).Value = $row.D
    $wsTag.Cells.Item($r, 6).Value = "int"
    $wsTag.Cells.Item($r, 7).Value = $row.D
    $wsTag.Cells.Item($r, 8).Value = "PCS"
    $wsTag.Cells.Item($r, 9).Value = 999999
    $wsTag.Cells.Item($r, 10).Value = "SUM"
    $wsTag.Cells.Item($r, 11).Value = "Active"
    $wsTag.Cells.Item($r, 12).Value = "SQL"
    $wsTag.Cells.Item($r, 13).Value = $row.M
    $wsTag.Cells.Item($r, 14).Value = $row.N
    $wsTag.Cells.Item($r, 15).Value = $row.O
    $wsTag.Cells.Item($r, 16).Value = $row.P
    $wsTag.Cells.Item($r, 17).Value = $row.Q
    $r = $r + 1
}

# Row 14: the original 168H-0012 row, shifted down from row 2, with a new
# tag_id (49) and site_id/asset_id columns appended
$wsTag.Cells.Item(14, 1).Value = 49
$wsTag.Cells.Item(14, 2).Value = "168H-0012.Cutter.Count"
$wsTag.Cells.Item(14, 3).Value = "168H-0012.Cutter.Count"
$wsTag.Cells.Item(14, 4).Value = "HA 12 Complete Counter"
$wsTag.Cells.Item(14, 5).Value = "NULL"
$wsTag.Cells.Item(14, 6).Value = "int"
$wsTag.Cells.Item(14, 7).Value = "NULL"
$wsTag.Cells.Item(14, 8).Value = "PCS"
$wsTag.Cells.Item(14, 9).Value = 999999
$wsTag.Cells.Item(14, 10).Value = "SUM"
$wsTag.Cells.Item(14, 11).Value = "Active"
$wsTag.Cells.Item(14, 12).Value = "SQL"
$wsTag.Cells.Item(14, 13).Value = "2019-12-11 15:29:17.567"
$wsTag.Cells.Item(14, 14).Value = "SQL manual entry"
$wsTag.Cells.Item(14, 15).Value = "2019-12-17 00:00:00.000"
$wsTag.Cells.Item(14, 16).Value = 119
$wsTag.Cells.Item(14, 17).Value = 112

$wsTag.Columns.Item(3).ColumnWidth = 9.140625
$wsTag.Range("N17").Select()

# ---------------------------------------------------------------------------
# 5. commonparameters sheet: site_id (col B) 120 -> 119
# ---------------------------------------------------------------------------
$wsCommon.Range("B2").Value = 119
$wsCommon.Range("J14").Select()

# ---------------------------------------------------------------------------
# 6. uom sheet: asset_id (col J) 120 -> 119
# ---------------------------------------------------------------------------
$wsUom = $wb.Worksheets.Item("uom")
$wsUom.Range("J2").Value = 119
$wsUom.Range("H17").Select()

# ---------------------------------------------------------------------------
# 7. unavailable sheet: selection only
# ---------------------------------------------------------------------------
$wsUnavail = $wb.Worksheets.Item("unavailable")
$wsUnavail.Range("F39").Select()

# ---------------------------------------------------------------------------
# 8. tfdusers sheet: asset_id (col G) 127 -> 119 for every data row, and this
#    sheet remains the active tab
# ---------------------------------------------------------------------------
$wsUsers = $wb.Worksheets.Item("tfdusers")
for ($r = 2; $r -le 42; $r++) {
    $wsUsers.Cells.Item($r, 7).Value = 119
}
$wsUsers.Activate()
$wsUsers.Range("L13").Select()
